# "M07 Froze Encoder 1"
# Re-run of the epoch-accuracy logging cell: the B column (per-epoch
# validation accuracy) is overwritten with the new run's numbers for the
# existing epochs (rows 2-101), the trailing "DisplayOutputs" rows
# (102-109) get their accuracy refreshed and their repr text updated to
# the new Python process's object memory address, and nine more
# freshly-logged rows are appended at the end (110-118, same repr text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column B: accuracy values for rows 2..118 (epochs 0..99 plus the
# trailing summary/extra rows) ----
$bData = @(
    0.9375, 0.875, 0.84375, 0.796875, 0.78125, 0.75, 0.734375, 0.703125,
    0.71875, 0.671875, 0.625, 0.640625, 0.609375, 0.59375, 0.546875,
    0.421875, 0.421875, 0.421875, 0.4375, 0.4375, 0.453125, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.421875, 0.421875, 0.421875, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375,
    0.5, 0.59375, 0.515625, 0.515625, 0.5, 0.5625, 0.5,
    0.453125, 0.578125, 0.453125, 0.546875, 0.5625, 0.5, 0.46875,
    0.546875, 0.3770491803278688
)

$bArr = New-Object 'object[,]' $bData.Length, 1
for ($i = 0; $i -lt $bData.Length; $i++) {
    $bArr[$i, 0] = $bData[$i]
}
$ws.Range("B2:B118").Value = $bArr

# ---- Column A: rows 102..118 hold the repr() of the DisplayOutputs
# object used for that training run; the memory address changed because
# this is a new Python process/run ----
$displayRepr = "<__main__.DisplayOutputs object at 0x7f19e8755cd0>"

$aArr = New-Object 'object[,]' 17, 1
for ($i = 0; $i -lt 17; $i++) {
    $aArr[$i, 0] = $displayRepr
}
$ws.Range("A102:A118").Value = $aArr

# ---- View state: the author selected the whole sheet (Ctrl+A) after
# scrolling back to the top, so the saved view no longer has a
# topLeftCell/F93 selection leftover from before ----
[void]$ws.Range("A1").Select()
[void]$ws.Cells.Select()
